# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at row 105 (pushing the existing
# rows 105-142 down to 106-143) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 105..142 down by one row.
$ws.Rows.Item(105).Insert()

# Fill in the new row 105 with the latest weekly data point.
$ws.Range("A105").Value = 9
$ws.Range("B105").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C105").Value = "Metropolitana"
$ws.Range("D105").Value = 44468
$ws.Range("E105").Value = 13
$ws.Range("F105").Value = 100112026
$ws.Range("G105").Value = "Haba"
$ws.Range("H105").Value = "Sin especificar"
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 79
$ws.Range("K105").Value = 8000
$ws.Range("L105").Value = 9000
$ws.Range("M105").Value = 8620
$ws.Range("N105").Value = "`$/saco 25 kilos"
$ws.Range("O105").Value = "Región Metropolitana"
$ws.Range("P105").Value = 345
$ws.Range("Q105").Value = 25
$ws.Range("R105").Value = "Hortaliza"
